# Weekly fruit/vegetable price update: insert one new daily record for
# "Ají" at Vega Monumental Concepción, shifting the existing rows 140-175
# down to 141-176.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data rows (140..175) down by one to make room for the
# new record; Excel's row insert shifts everything below automatically.
$ws.Rows("140:140").Insert()

# Populate the newly-inserted row 140 with the new price record.
$ws.Range("A140").Value = 11
$ws.Range("B140").Value = "Vega Monumental Concepción"
$ws.Range("C140").Value = "Bíobío"
$ws.Range("D140").Value = 44995
$ws.Range("E140").Value = 8
$ws.Range("F140").Value = 100112021
$ws.Range("G140").Value = "Ají"
$ws.Range("H140").Value = "Americana (o)"
$ws.Range("I140").Value = "Primera"
$ws.Range("J140").Value = 100
$ws.Range("K140").Value = 20000
$ws.Range("L140").Value = 22000
$ws.Range("M140").Value = 21000
$ws.Range("N140").Value = "$/saco 25 kilos"
$ws.Range("O140").Value = "Región Metropolitana"
$ws.Range("P140").Value = 840
$ws.Range("Q140").Value = 25
$ws.Range("R140").Value = "Hortaliza"
